$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "244.41"
Set-TextValue 2 5 "-0.89%"
Set-TextValue 3 4 "27.61"
Set-TextValue 3 5 "6.26%"
Set-TextValue 4 4 "5.129"
Set-TextValue 4 5 "0.72%"
Set-TextValue 5 5 "1.66%"
Set-TextValue 6 4 "6.495"
Set-TextValue 6 5 "0.26%"
Set-TextValue 7 4 "0.8204"
Set-TextValue 7 5 "0.83%"
Set-TextValue 8 4 "0.8565"
Set-TextValue 8 5 "1.15%"
Set-TextValue 9 4 "0.06945"
Set-TextValue 9 5 "0.08%"
Set-TextValue 10 4 "0.02879"
Set-TextValue 10 5 "2.17%"
Set-TextValue 11 4 "0.09397"
Set-TextValue 11 5 "0.10%"
Set-TextValue 12 4 "0.001524"
Set-TextValue 12 5 "0.84%"
Set-TextValue 13 4 "0.04090"
Set-TextValue 13 5 "-11.83%"
Set-TextValue 14 4 "0.0006030"
Set-TextValue 14 5 "1.17%"
Set-TextValue 15 4 "0.006214"
Set-TextValue 15 5 "-0.59%"
Set-TextValue 16 5 "-2.40%"
Set-TextValue 17 4 "3.012"
Set-TextValue 17 5 "-0.39%"
Set-TextValue 18 4 "2.185"
Set-TextValue 18 5 "6.29%"
Set-TextValue 19 4 "0.3151"
Set-TextValue 19 5 "1.26%"
Set-TextValue 20 4 "0.1334"
Set-TextValue 20 5 "0.16%"
Set-TextValue 21 4 "0.03224"
Set-TextValue 21 5 "1.75%"
Set-TextValue 22 5 "-1.53%"
Set-TextValue 23 4 "3.553"
Set-TextValue 23 5 "-5.21%"
Set-TextValue 24 4 "0.1374"
Set-TextValue 24 5 "-0.05%"
Set-TextValue 25 5 "-2.65%"
Set-TextValue 26 4 "0.004475"
Set-TextValue 26 5 "-1.68%"
Set-TextValue 27 4 "0.0001180"
Set-TextValue 27 5 "22.85%"
Set-TextValue 28 5 "-27.44%"
Set-TextValue 40 4 "0.03720"
Set-TextValue 40 5 "1.82%"
Set-TextValue 41 4 "0.005918"
Set-TextValue 41 5 "-3.79%"
Set-TextValue 42 4 "0.1058"
Set-TextValue 42 5 "0.46%"
Set-TextValue 43 4 "0.002399"
Set-TextValue 43 5 "-7.33%"
Set-TextValue 44 4 "0.009593"
Set-TextValue 44 5 "20.30%"
Set-TextValue 45 4 "0.00005113"
Set-TextValue 45 5 "-5.24%"
Set-TextValue 46 5 "-0.03%"
Set-TextValue 47 4 "0.1010"
Set-TextValue 47 5 "-30.35%"
Set-TextValue 48 4 "0.002568"
Set-TextValue 48 5 "7.05%"
Set-TextValue 49 4 "0.00002100"
Set-TextValue 49 5 "-0.03%"
Set-TextValue 50 4 "0.0002000"
Set-TextValue 50 5 "-0.03%"
